$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (M) mirroring the existing "2020" column (L):
# copy L3:L7 formatting into M3:M7, then set the new year's values.

$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)

$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)

$ws.Range("M4").Value = 2021
$ws.Range("M5").Value = 98
$ws.Range("M6").Value = 97
$ws.Range("M7").Value = 96

# Reset the view: drop the scrolled topLeftCell and the stray N13 selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
